$wb = $excel.ActiveWorkbook

# Sheet: ALC (sheet1.xml)
$ws = $wb.Worksheets.Item(1)
$ws.Range("H86").Value = 288352
$ws.Range("I86").Value = 3791
$ws.Range("J86").Value = 402176.4
$ws.Range("K86").Value = 3791
$ws.Range("L86").Value = 402176.4
$ws.Range("M86").Value = -2668
$ws.Range("N86").Value = -404422.4
$ws.Range("H89").Value = 288352
$ws.Range("I89").Value = 3791
$ws.Range("J89").Value = 402176.4
$ws.Range("K89").Value = 18955
$ws.Range("L89").Value = 2010882
$ws.Range("M89").Value = -13339
$ws.Range("N89").Value = -2022114
$ws.Range("H100").Value = 1153.1538
$ws.Range("I100").Value = 1088.6
$ws.Range("K100").Value = 1088.6
$ws.Range("M100").Value = -547.5999999999999
$ws.Range("H135").Value = 1577.8462
$ws.Range("I135").Value = 539.6667
$ws.Range("K135").Value = 4857.0003
$ws.Range("M135").Value = -2322.0003
$ws.Range("H137").Value = 1307.6666
$ws.Range("I137").Value = 1239.3636
$ws.Range("K137").Value = 3718.0908
$ws.Range("M137").Value = -1168.0908
$ws.Range("H138").Value = 4788.0645
$ws.Range("I138").Value = 3839.8333
$ws.Range("J138").Value = 5015.64
$ws.Range("K138").Value = 11519.4999
$ws.Range("L138").Value = 15046.92
$ws.Range("M138").Value = -6379.499899999999
$ws.Range("N138").Value = -25326.92

# Sheet: ARM (sheet2.xml)
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 3211198
$ws.Range("I32").Value = 3046467.5
$ws.Range("K32").Value = 3046467.5
$ws.Range("M32").Value = -3046180.5
$ws.Range("H45").Value = 944.4
$ws.Range("I45").Value = 902.6667
$ws.Range("J45").Value = 1007
$ws.Range("K45").Value = 902.6667
$ws.Range("L45").Value = 1007
$ws.Range("M45").Value = -525.6667
$ws.Range("N45").Value = -1761
$ws.Range("H63").Value = 3981.182
$ws.Range("I63").Value = 4029.3
$ws.Range("K63").Value = 4029.3
$ws.Range("M63").Value = -3343.3
$ws.Range("H66").Value = 3981.182
$ws.Range("I66").Value = 4029.3
$ws.Range("K66").Value = 20146.5
$ws.Range("M66").Value = -16714.5
$ws.Range("H74").Value = 1070.5
$ws.Range("I74").Value = 1330
$ws.Range("K74").Value = 1330
$ws.Range("M74").Value = -456
$ws.Range("H77").Value = 1070.5
$ws.Range("I77").Value = 1330
$ws.Range("K77").Value = 6650
$ws.Range("M77").Value = -2282
$ws.Range("H88").Value = 2450.25
$ws.Range("J88").Value = 2299.4285
$ws.Range("L88").Value = 2299.4285
$ws.Range("N88").Value = -3111.4285
$ws.Range("H91").Value = 2450.25
$ws.Range("J91").Value = 2299.4285
$ws.Range("L91").Value = 2299.4285
$ws.Range("N91").Value = -5107.4285
$ws.Range("H113").Value = 129849
$ws.Range("J113").Value = 129849
$ws.Range("L113").Value = 129849
$ws.Range("N113").Value = -138527
$ws.Range("H138").Value = 648342.2
$ws.Range("J138").Value = 648342.2
$ws.Range("L138").Value = 648342.2
$ws.Range("N138").Value = -658622.2

# Sheet: BSM (sheet3.xml)
$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 1627.3
$ws.Range("I86").Value = 1627.3
$ws.Range("K86").Value = 1627.3
$ws.Range("M86").Value = -504.3
$ws.Range("H89").Value = 1627.3
$ws.Range("I89").Value = 1627.3
$ws.Range("K89").Value = 8136.5
$ws.Range("M89").Value = -2520.5
$ws.Range("H99").Value = 1957.6154
$ws.Range("I99").Value = 1970.75
$ws.Range("J99").Value = 1800
$ws.Range("K99").Value = 1970.75
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = -472.75
$ws.Range("N99").Value = -4796
$ws.Range("H105").Value = 1848.5
$ws.Range("I105").Value = 1848.5
$ws.Range("K105").Value = 1848.5
$ws.Range("M105").Value = -101.5
$ws.Range("H107").Value = 498.7143
$ws.Range("I107").Value = 497.75
$ws.Range("K107").Value = 497.75
$ws.Range("M107").Value = 1422.25

# Sheet: CRP (sheet4.xml)
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = 183.14285
$ws.Range("I22").Value = 230
$ws.Range("K22").Value = 230
$ws.Range("M22").Value = 120
$ws.Range("H62").Value = 999
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 999
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 999
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -2247
$ws.Range("H65").Value = 999
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 999
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 4995
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -11235
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# Sheet: CUL (sheet5.xml)
$ws = $wb.Worksheets.Item(5)
$ws.Range("H12").Value = 183.3
$ws.Range("J12").Value = 149.4
$ws.Range("L12").Value = 448.2
$ws.Range("N12").Value = -794.2
$ws.Range("H33").Value = 966.3333
$ws.Range("J33").Value = 966.3333
$ws.Range("L33").Value = 5797.9998
$ws.Range("N33").Value = -6363.9998
$ws.Range("H34").Value = 999.6667
$ws.Range("I34").Value = 499
$ws.Range("J34").Value = 1250
$ws.Range("K34").Value = 1497
$ws.Range("L34").Value = 3750
$ws.Range("M34").Value = -1413
$ws.Range("N34").Value = -3918
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H44").Value = 507.64706
$ws.Range("I44").Value = 843.75
$ws.Range("J44").Value = 404.23077
$ws.Range("K44").Value = 2531.25
$ws.Range("L44").Value = 1212.69231
$ws.Range("M44").Value = -2133.25
$ws.Range("N44").Value = -2008.69231
$ws.Range("H51").Value = 1004
$ws.Range("I51").Value = 1004
$ws.Range("K51").Value = 3012
$ws.Range("M51").Value = -2552
$ws.Range("H55").Value = 5992.5
$ws.Range("J55").Value = 5992.5
$ws.Range("L55").Value = 17977.5
$ws.Range("N55").Value = -18331.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 900
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2700
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -4322
$ws.Range("H70").Value = 13856.429
$ws.Range("I70").Value = 3498.5
$ws.Range("J70").Value = 17999.6
$ws.Range("K70").Value = 10495.5
$ws.Range("L70").Value = 53998.8
$ws.Range("M70").Value = -10180.5
$ws.Range("N70").Value = -54628.8
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 900
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 8100
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -16212
$ws.Range("H73").Value = 13856.429
$ws.Range("I73").Value = 3498.5
$ws.Range("J73").Value = 17999.6
$ws.Range("K73").Value = 10495.5
$ws.Range("L73").Value = 53998.8
$ws.Range("M73").Value = -9403.5
$ws.Range("N73").Value = -56182.8

# Sheet: GSM (sheet6.xml)
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 4064.6667
$ws.Range("I70").Value = 3997
$ws.Range("K70").Value = 3997
$ws.Range("M70").Value = -3727
$ws.Range("H73").Value = 4064.6667
$ws.Range("I73").Value = 3997
$ws.Range("K73").Value = 3997
$ws.Range("M73").Value = -3061
$ws.Range("H80").Value = 3116.25
$ws.Range("I80").Value = 2468.3333
$ws.Range("J80").Value = 3505
$ws.Range("K80").Value = 2468.3333
$ws.Range("L80").Value = 3505
$ws.Range("M80").Value = -1470.3333
$ws.Range("N80").Value = -5501
$ws.Range("H83").Value = 3116.25
$ws.Range("I83").Value = 2468.3333
$ws.Range("J83").Value = 3505
$ws.Range("K83").Value = 12341.6665
$ws.Range("L83").Value = 17525
$ws.Range("M83").Value = -7349.666499999999
$ws.Range("N83").Value = -27509
$ws.Range("H132").Value = 8971.883
$ws.Range("I132").Value = 8971.883
$ws.Range("K132").Value = 26915.649
$ws.Range("M132").Value = -24385.649

# Sheet: LTW (sheet7.xml)
$ws = $wb.Worksheets.Item(7)
$ws.Range("H61").Value = 480
$ws.Range("I61").Value = 480
$ws.Range("K61").Value = 480
$ws.Range("M61").Value = -278
$ws.Range("H106").Value = 23928.5
$ws.Range("J106").Value = 23928.5
$ws.Range("L106").Value = 23928.5
$ws.Range("N106").Value = -26452.5
$ws.Range("H113").Value = 480
$ws.Range("I113").Value = 480
$ws.Range("K113").Value = 480
$ws.Range("M113").Value = 1690
$ws.Range("H122").Value = 5538.0605
$ws.Range("I122").Value = 4617.8945
$ws.Range("J122").Value = 6786.857
$ws.Range("K122").Value = 13853.6835
$ws.Range("L122").Value = 20360.571
$ws.Range("M122").Value = -11403.6835
$ws.Range("N122").Value = -25260.571
$ws.Range("H127").Value = 49999.5
$ws.Range("J127").Value = 49999.5
$ws.Range("L127").Value = 49999.5
$ws.Range("N127").Value = -59919.5
$ws.Range("H136").Value = 3858.8
$ws.Range("I136").Value = 3858.8
$ws.Range("K136").Value = 11576.4
$ws.Range("M136").Value = -9026.400000000001

# Sheet: WVR (sheet8.xml)
$ws = $wb.Worksheets.Item(8)
$ws.Range("H41").Value = 19610.75
$ws.Range("I41").Value = 19563.666
$ws.Range("K41").Value = 19563.666
$ws.Range("M41").Value = -19173.666
